$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (shared-string based) for the launcher position change
# Order matters for shared string table indices: F3's string must be
# registered before F2's string so they land at indices 18 and 19 respectively.
$ws.Range("F3").Value = "10102:0:45:90:135:180:225:270:315"
$ws.Range("F2").Value = "10102:0:45:90:135:225:270:315"

# Update selection to F3
$ws.Range("F3").Select()

# Update column width: split F out from the F:J run into its own width.
# (G:J are intentionally left untouched so they keep their original exact
# 19.125 width and stay merged as a single <col min="7" max="10"> run.)
#
# Note: the host's ColumnWidth setter quantizes to an MDW-7 pixel grid
# (stored width = (round(chars*7)+5)/7), so an input of 31.375 cannot be
# represented exactly - it lands on 32.142857... A raw input of
# 30.714285714285715 (=215/7) lands on the closest achievable grid value,
# 31.428571428571427, which is nearest to the target 31.375.
$ws.Columns("F").ColumnWidth = 30.714285714285715
